$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# The document starts as 31 separate paragraphs ("1".."31"), each with
# identical paragraph formatting:
#   <w:spacing w:line="360" w:lineRule="auto" w:after="0"/><w:jc w:val="left"/>
#
# The target is a single paragraph whose pPr keeps only the line spacing
# (no w:after, no w:jc), containing one run with two sentences of Lorem
# ipsum text separated by a manual line break (<w:br/>).
# ---------------------------------------------------------------------

$loremPart1 = "Lorem ipsum odor amet, consectetuer adipiscing elit. Congue vel parturient sapien volutpat porttitor malesuada mus. Volutpat sociosqu nisi cubilia himenaeos sed in nisl leo. Dis venenatis ullamcorper pharetra; penatibus blandit arcu justo dignissim nullam. Dolor a sodales, nostra lacinia nascetur faucibus. Sodales volutpat mattis suscipit morbi"
$loremPart2 = "elementum sapien convallis nec egestas. Dignissim lacinia dolor placerat nulla porta natoque fames, sem non. Venenatis facilisi dapibus pellentesque netus etiam id blandit. Porttitor integer nec urna posuere rhoncus faucibus."

# Step 1: Trim the first paragraph's own properties down to just the
# spacing (drop w:after="0" and w:jc="left") while it still has a
# following paragraph -- inserting a <w:p> fragment at the very end of
# the document's content instead would split off a stray empty
# paragraph, so this has to happen before the other 30 paragraphs are
# collapsed away.
$p1 = $d.Paragraphs.First
$pMark = $d.Range($p1.Range.End - 1, $p1.Range.End)
$pprXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/></w:pPr><w:r><w:t>1</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pMark.InsertXML($pprXml)

# Step 2: Merge all 31 paragraphs into one by deleting every paragraph
# mark (^p) in the document -- this concatenates the runs "1".."31"
# into the (now reformatted) first paragraph and leaves just one
# paragraph behind.
$d.Content.Find.Execute("^p", $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# Step 3: Replace the merged run text with the Lorem ipsum content,
# keeping a manual line break (<w:br/>) between the two sentences. The
# replacement range stops one character short of the paragraph's end so
# only the run content is swapped, leaving the paragraph mark (and the
# pPr fixed in step 1) untouched.
$p1 = $d.Paragraphs.First
$body = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$runXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>' + $loremPart1 + '</w:t><w:br/><w:t>' + $loremPart2 + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$body.InsertXML($runXml)

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
Write-Output ("Text=" + $d.Paragraphs.First.Range.Text)
